$d = $word.ActiveDocument
$d.Content.Find.Execute("49+19=68", $true, $false, $false, $false, $false, $true, 1, $false, "86+7=93", 2) | Out-Null
$d.Content.Find.Execute("97-62=35", $true, $false, $false, $false, $false, $true, 1, $false, "42-40=2", 2) | Out-Null
$d.Content.Find.Execute("1+85=86", $true, $false, $false, $false, $false, $true, 1, $false, "40-27=13", 2) | Out-Null
$d.Content.Find.Execute("35-32=3", $true, $false, $false, $false, $false, $true, 1, $false, "5+14=19", 2) | Out-Null
$d.Content.Find.Execute("74-11=63", $true, $false, $false, $false, $false, $true, 1, $false, "95-13=82", 2) | Out-Null
$d.Content.Find.Execute("66-17=49", $true, $false, $false, $false, $false, $true, 1, $false, "18+72=90", 2) | Out-Null
$d.Content.Find.Execute("81-66=15", $true, $false, $false, $false, $false, $true, 1, $false, "80-8=72", 2) | Out-Null
$d.Content.Find.Execute("51-30=21", $true, $false, $false, $false, $false, $true, 1, $false, "46-21=25", 2) | Out-Null
$d.Content.Find.Execute("56-42=14", $true, $false, $false, $false, $false, $true, 1, $false, "71+9=80", 2) | Out-Null
$d.Content.Find.Execute("57-51=6", $true, $false, $false, $false, $false, $true, 1, $false, "13-5=8", 2) | Out-Null
$d.Content.Find.Execute("41+55=96", $true, $false, $false, $false, $false, $true, 1, $false, "22-0=22", 2) | Out-Null
$d.Content.Find.Execute("5+89=94", $true, $false, $false, $false, $false, $true, 1, $false, "62-13=49", 2) | Out-Null
$d.Content.Find.Execute("50+11=61", $true, $false, $false, $false, $false, $true, 1, $false, "73-8=65", 2) | Out-Null
$d.Content.Find.Execute("1+19=20", $true, $false, $false, $false, $false, $true, 1, $false, "2+52=54", 2) | Out-Null
$d.Content.Find.Execute("85+7=92", $true, $false, $false, $false, $false, $true, 1, $false, "44+45=89", 2) | Out-Null
$d.Content.Find.Execute("57+0=57", $true, $false, $false, $false, $false, $true, 1, $false, "13-1=12", 2) | Out-Null
$d.Content.Find.Execute("19+32=51", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=84", 2) | Out-Null
$d.Content.Find.Execute("40+28=68", $true, $false, $false, $false, $false, $true, 1, $false, "88-46=42", 2) | Out-Null
$d.Content.Find.Execute("86-31=55", $true, $false, $false, $false, $false, $true, 1, $false, "60+15=75", 2) | Out-Null
$d.Content.Find.Execute("65+21=86", $true, $false, $false, $false, $false, $true, 1, $false, "70-33=37", 2) | Out-Null
$d.Content.Find.Execute("44+8=52", $true, $false, $false, $false, $false, $true, 1, $false, "14+4=18", 2) | Out-Null
$d.Content.Find.Execute("57+33=90", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=29", 2) | Out-Null
$d.Content.Find.Execute("66-7=59", $true, $false, $false, $false, $false, $true, 1, $false, "49-48=1", 2) | Out-Null
$d.Content.Find.Execute("54-19=35", $true, $false, $false, $false, $false, $true, 1, $false, "16+18=34", 2) | Out-Null
$d.Content.Find.Execute("84-47=37", $true, $false, $false, $false, $false, $true, 1, $false, "93-39=54", 2) | Out-Null
$d.Content.Find.Execute("89+3=92", $true, $false, $false, $false, $false, $true, 1, $false, "75-68=7", 2) | Out-Null
$d.Content.Find.Execute("49-8=41", $true, $false, $false, $false, $false, $true, 1, $false, "67-59=8", 2) | Out-Null
$d.Content.Find.Execute("46-20=26", $true, $false, $false, $false, $false, $true, 1, $false, "37-25=12", 2) | Out-Null
$d.Content.Find.Execute("59+25=84", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=26", 2) | Out-Null
$d.Content.Find.Execute("9+89=98", $true, $false, $false, $false, $false, $true, 1, $false, "33+60=93", 2) | Out-Null
$d.Content.Find.Execute("36-8=28", $true, $false, $false, $false, $false, $true, 1, $false, "79-78=1", 2) | Out-Null
$d.Content.Find.Execute("93-70=23", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=46", 2) | Out-Null
$d.Content.Find.Execute("76+19=95", $true, $false, $false, $false, $false, $true, 1, $false, "3+20=23", 2) | Out-Null
$d.Content.Find.Execute("56+19=75", $true, $false, $false, $false, $false, $true, 1, $false, "19+0=19", 2) | Out-Null
$d.Content.Find.Execute("15+6=21", $true, $false, $false, $false, $false, $true, 1, $false, "74-39=35", 2) | Out-Null
$d.Content.Find.Execute("55-48=7", $true, $false, $false, $false, $false, $true, 1, $false, "72-59=13", 2) | Out-Null
$d.Content.Find.Execute("20+8=28", $true, $false, $false, $false, $false, $true, 1, $false, "91-62=29", 2) | Out-Null
$d.Content.Find.Execute("91-15=76", $true, $false, $false, $false, $false, $true, 1, $false, "45+22=67", 2) | Out-Null
$d.Content.Find.Execute("88-38=50", $true, $false, $false, $false, $false, $true, 1, $false, "34+53=87", 2) | Out-Null
$d.Content.Find.Execute("34-25=9", $true, $false, $false, $false, $false, $true, 1, $false, "47+24=71", 2) | Out-Null
$d.Content.Find.Execute("83-75=8", $true, $false, $false, $false, $false, $true, 1, $false, "73-70=3", 2) | Out-Null
$d.Content.Find.Execute("22+52=74", $true, $false, $false, $false, $false, $true, 1, $false, "25-0=25", 2) | Out-Null
$d.Content.Find.Execute("97-32=65", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=51", 2) | Out-Null
$d.Content.Find.Execute("3+86=89", $true, $false, $false, $false, $false, $true, 1, $false, "72-38=34", 2) | Out-Null
$d.Content.Find.Execute("81-69=12", $true, $false, $false, $false, $false, $true, 1, $false, "8+59=67", 2) | Out-Null
$d.Content.Find.Execute("92-41=51", $true, $false, $false, $false, $false, $true, 1, $false, "44+55=99", 2) | Out-Null
$d.Content.Find.Execute("40+42=82", $true, $false, $false, $false, $false, $true, 1, $false, "34+13=47", 2) | Out-Null
$d.Content.Find.Execute("5+79=84", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("87-70=17", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=35", 2) | Out-Null
$d.Content.Find.Execute("70-1=69", $true, $false, $false, $false, $false, $true, 1, $false, "50-28=22", 2) | Out-Null
$d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, 1, $false, "29-27=2", 2) | Out-Null
$d.Content.Find.Execute("2+41=43", $true, $false, $false, $false, $false, $true, 1, $false, "8+38=46", 2) | Out-Null
$d.Content.Find.Execute("91-11=80", $true, $false, $false, $false, $false, $true, 1, $false, "21+45=66", 2) | Out-Null
$d.Content.Find.Execute("49-30=19", $true, $false, $false, $false, $false, $true, 1, $false, "80-18=62", 2) | Out-Null
$d.Content.Find.Execute("43+44=87", $true, $false, $false, $false, $false, $true, 1, $false, "28+27=55", 2) | Out-Null
$d.Content.Find.Execute("73-52=21", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("87-59=28", $true, $false, $false, $false, $false, $true, 1, $false, "79-31=48", 2) | Out-Null
$d.Content.Find.Execute("8+67=75", $true, $false, $false, $false, $false, $true, 1, $false, "45-0=45", 2) | Out-Null
$d.Content.Find.Execute("4+1=5", $true, $false, $false, $false, $false, $true, 1, $false, "3+34=37", 2) | Out-Null
$d.Content.Find.Execute("65-3=62", $true, $false, $false, $false, $false, $true, 1, $false, "33+8=41", 2) | Out-Null
$d.Content.Find.Execute("22+11=33", $true, $false, $false, $false, $false, $true, 1, $false, "64+22=86", 2) | Out-Null
$d.Content.Find.Execute("25+39=64", $true, $false, $false, $false, $false, $true, 1, $false, "71-29=42", 2) | Out-Null
$d.Content.Find.Execute("77+18=95", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=45", 2) | Out-Null
$d.Content.Find.Execute("17+32=49", $true, $false, $false, $false, $false, $true, 1, $false, "34+22=56", 2) | Out-Null
$d.Content.Find.Execute("12+62=74", $true, $false, $false, $false, $false, $true, 1, $false, "3+3=6", 2) | Out-Null
$d.Content.Find.Execute("66-44=22", $true, $false, $false, $false, $false, $true, 1, $false, "29+12=41", 2) | Out-Null
$d.Content.Find.Execute("68-36=32", $true, $false, $false, $false, $false, $true, 1, $false, "8+11=19", 2) | Out-Null
$d.Content.Find.Execute("43-41=2", $true, $false, $false, $false, $false, $true, 1, $false, "14+50=64", 2) | Out-Null
$d.Content.Find.Execute("75-24=51", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=55", 2) | Out-Null
$d.Content.Find.Execute("41+8=49", $true, $false, $false, $false, $false, $true, 1, $false, "61-31=30", 2) | Out-Null
$d.Content.Find.Execute("70+3=73", $true, $false, $false, $false, $false, $true, 1, $false, "8+19=27", 2) | Out-Null
$d.Content.Find.Execute("62+7=69", $true, $false, $false, $false, $false, $true, 1, $false, "45+29=74", 2) | Out-Null
$d.Content.Find.Execute("86+5=91", $true, $false, $false, $false, $false, $true, 1, $false, "21+15=36", 2) | Out-Null
$d.Content.Find.Execute("72+16=88", $true, $false, $false, $false, $false, $true, 1, $false, "26+19=45", 2) | Out-Null
$d.Content.Find.Execute("45-20=25", $true, $false, $false, $false, $false, $true, 1, $false, "19+31=50", 2) | Out-Null
$d.Content.Find.Execute("89-74=15", $true, $false, $false, $false, $false, $true, 1, $false, "2+15=17", 2) | Out-Null
$d.Content.Find.Execute("38+26=64", $true, $false, $false, $false, $false, $true, 1, $false, "66-41=25", 2) | Out-Null
$d.Content.Find.Execute("14+26=40", $true, $false, $false, $false, $false, $true, 1, $false, "66+23=89", 2) | Out-Null
$d.Content.Find.Execute("83-72=11", $true, $false, $false, $false, $false, $true, 1, $false, "84+3=87", 2) | Out-Null
$d.Content.Find.Execute("97-51=46", $true, $false, $false, $false, $false, $true, 1, $false, "87-49=38", 2) | Out-Null
$d.Content.Find.Execute("93-42=51", $true, $false, $false, $false, $false, $true, 1, $false, "9+41=50", 2) | Out-Null
$d.Content.Find.Execute("93-88=5", $true, $false, $false, $false, $false, $true, 1, $false, "12+72=84", 2) | Out-Null
$d.Content.Find.Execute("77-23=54", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=28", 2) | Out-Null
$d.Content.Find.Execute("3+2=5", $true, $false, $false, $false, $false, $true, 1, $false, "42-24=18", 2) | Out-Null
$d.Content.Find.Execute("79-59=20", $true, $false, $false, $false, $false, $true, 1, $false, "1+94=95", 2) | Out-Null
$d.Content.Find.Execute("11+20=31", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("44+41=85", $true, $false, $false, $false, $false, $true, 1, $false, "93-0=93", 2) | Out-Null
$d.Content.Find.Execute("13-3=10", $true, $false, $false, $false, $false, $true, 1, $false, "52+5=57", 2) | Out-Null
$d.Content.Find.Execute("99-17=82", $true, $false, $false, $false, $false, $true, 1, $false, "96-26=70", 2) | Out-Null
$d.Content.Find.Execute("76+8=84", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=78", 2) | Out-Null
$d.Content.Find.Execute("67-61=6", $true, $false, $false, $false, $false, $true, 1, $false, "13+29=42", 2) | Out-Null
$d.Content.Find.Execute("38+25=63", $true, $false, $false, $false, $false, $true, 1, $false, "7+36=43", 2) | Out-Null
$d.Content.Find.Execute("65-61=4", $true, $false, $false, $false, $false, $true, 1, $false, "22+2=24", 2) | Out-Null
$d.Content.Find.Execute("89-14=75", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("80-77=3", $true, $false, $false, $false, $false, $true, 1, $false, "46-33=13", 2) | Out-Null
$d.Content.Find.Execute("15-9=6", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("98-91=7", $true, $false, $false, $false, $false, $true, 1, $false, "19-8=11", 2) | Out-Null
$d.Content.Find.Execute("45+33=78", $true, $false, $false, $false, $false, $true, 1, $false, "23+47=70", 2) | Out-Null
$d.Content.Find.Execute("34+7=41", $true, $false, $false, $false, $false, $true, 1, $false, "81-44=37", 2) | Out-Null
$d.Content.Find.Execute("84-19=65", $true, $false, $false, $false, $false, $true, 1, $false, "47+35=82", 2) | Out-Null
